$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header cell C1 from "Task Name" to "Name"
$ws.Range("C1").Value = "Name"

# Update the selected/active cell from C2 to C1
$ws.Range("C1").Select()
